# Weekly refresh of the price series: a new day's reading is inserted at
# the top of the data block (row 18), pushing the existing rows down by
# one and re-appending the last existing row's data at the new bottom
# (row 37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 18:36 down to 19:37, duplicating row 36's old values into the
# new row 37 and leaving a blank row 18 to populate with the new reading.
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 44669
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 100112052
$ws.Range("G18").Value = "Albahaca"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 950
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 975
$ws.Range("N18").Value = "$/paquete"
$ws.Range("O18").Value = "Región de Arica y Parinacota"
$ws.Range("P18").Value = 975
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"
